# Recompute ligand/receptor expression & specificity values for the
# Cadm1-Cadm1 LR-pair sheet using the updated TPM normalization.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 0.82446933333333339
$ws.Range("H2").Value = 2.47340800000000005
$ws.Range("I2").Value = 0.045311626617283135
$ws.Range("J2").Value = 0.04531162661728313
$ws.Range("M2").Value = 0.82446933333333339
$ws.Range("N2").Value = 2.47340800000000005
$ws.Range("O2").Value = 0.045311626617283135
$ws.Range("P2").Value = 0.04531162661728313
$ws.Range("Q2").Value = 0.67974968160711124
$ws.Range("R2").Value = 6.11774713446400042
$ws.Range("S2").Value = 0.0020531435067040814
$ws.Range("T2").Value = 0.002053143506704081

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 0.82446933333333339
$ws.Range("H3").Value = 2.47340800000000005
$ws.Range("I3").Value = 0.045311626617283135
$ws.Range("J3").Value = 0.04531162661728313
$ws.Range("M3").Value = 0.01719833333333333
$ws.Range("N3").Value = 0.051595
$ws.Range("O3").Value = 0.0009451952024569835
$ws.Range("P3").Value = 0.0009451952024569837
$ws.Range("Q3").Value = 0.014179498417777775
$ws.Range("R3").Value = 0.12761548576000001
$ws.Range("S3").Value = 0.00004282833209417818
$ws.Range("T3").Value = 0.00004282833209417818

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 0.82446933333333339
$ws.Range("H4").Value = 2.47340800000000005
$ws.Range("I4").Value = 0.045311626617283135
$ws.Range("J4").Value = 0.04531162661728313
$ws.Range("M4").Value = 3.18231699999999984
$ws.Range("N4").Value = 9.54695099999999996
$ws.Range("O4").Value = 0.17489547985835649
$ws.Range("P4").Value = 0.17489547985835646
$ws.Range("Q4").Value = 2.62372277544533317
$ws.Range("R4").Value = 23.61350497900799894
$ws.Range("S4").Value = 0.007924798680392413
$ws.Range("T4").Value = 0.00792479868039241

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 0.82446933333333339
$ws.Range("H5").Value = 2.47340800000000005
$ws.Range("I5").Value = 0.045311626617283135
$ws.Range("J5").Value = 0.04531162661728313
$ws.Range("M5").Value = 14.17155133333332984
$ws.Range("N5").Value = 42.51465400000000017
$ws.Range("O5").Value = 0.77884769832190337
$ws.Range("P5").Value = 0.77884769832190348
$ws.Range("Q5").Value = 11.68400948009244189
$ws.Range("R5").Value = 105.15608532083200544
$ws.Range("S5").Value = 0.03529085609809246
$ws.Range("T5").Value = 0.03529085609809246

# Row 6: FAPs -> ECs
$ws.Range("G6").Value = 0.01719833333333333
$ws.Range("H6").Value = 0.051595
$ws.Range("I6").Value = 0.0009451952024569835
$ws.Range("J6").Value = 0.0009451952024569837
$ws.Range("M6").Value = 0.82446933333333339
$ws.Range("N6").Value = 2.47340800000000005
$ws.Range("O6").Value = 0.045311626617283135
$ws.Range("P6").Value = 0.04531162661728313
$ws.Range("Q6").Value = 0.014179498417777775
$ws.Range("R6").Value = 0.12761548576000001
$ws.Range("S6").Value = 0.00004282833209417818
$ws.Range("T6").Value = 0.00004282833209417818

# Row 7: FAPs -> FAPs
$ws.Range("G7").Value = 0.01719833333333333
$ws.Range("H7").Value = 0.051595
$ws.Range("I7").Value = 0.0009451952024569835
$ws.Range("J7").Value = 0.0009451952024569837
$ws.Range("M7").Value = 0.01719833333333333
$ws.Range("N7").Value = 0.051595
$ws.Range("O7").Value = 0.0009451952024569835
$ws.Range("P7").Value = 0.0009451952024569837
$ws.Range("Q7").Value = 0.0002957826694444443
$ws.Range("R7").Value = 0.002662044025
$ws.Range("S7").Value = 0.000000893393970747698
$ws.Range("T7").Value = 0.0000008933939707476985

# Row 8: FAPs -> MuSCs
$ws.Range("G8").Value = 0.01719833333333333
$ws.Range("H8").Value = 0.051595
$ws.Range("I8").Value = 0.0009451952024569835
$ws.Range("J8").Value = 0.0009451952024569837
$ws.Range("M8").Value = 3.18231699999999984
$ws.Range("N8").Value = 9.54695099999999996
$ws.Range("O8").Value = 0.17489547985835649
$ws.Range("P8").Value = 0.17489547985835646
$ws.Range("Q8").Value = 0.05473054853833332
$ws.Range("R8").Value = 0.49257493684500003
$ws.Range("S8").Value = 0.00016531036849353053
$ws.Range("T8").Value = 0.00016531036849353056

# Row 9: FAPs -> Resolving-Mac
$ws.Range("G9").Value = 0.01719833333333333
$ws.Range("H9").Value = 0.051595
$ws.Range("I9").Value = 0.0009451952024569835
$ws.Range("J9").Value = 0.0009451952024569837
$ws.Range("M9").Value = 14.17155133333332984
$ws.Range("N9").Value = 42.51465400000000017
$ws.Range("O9").Value = 0.77884769832190337
$ws.Range("P9").Value = 0.77884769832190348
$ws.Range("Q9").Value = 0.24372706368111099
$ws.Range("R9").Value = 2.19354357312999992
$ws.Range("S9").Value = 0.0007361631078985271
$ws.Range("T9").Value = 0.0007361631078985273

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 3.18231699999999984
$ws.Range("H10").Value = 9.54695099999999996
$ws.Range("I10").Value = 0.17489547985835649
$ws.Range("J10").Value = 0.17489547985835646
$ws.Range("M10").Value = 0.82446933333333339
$ws.Range("N10").Value = 2.47340800000000005
$ws.Range("O10").Value = 0.045311626617283135
$ws.Range("P10").Value = 0.04531162661728313
$ws.Range("Q10").Value = 2.62372277544533317
$ws.Range("R10").Value = 23.61350497900799894
$ws.Range("S10").Value = 0.007924798680392413
$ws.Range("T10").Value = 0.00792479868039241

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 3.18231699999999984
$ws.Range("H11").Value = 9.54695099999999996
$ws.Range("I11").Value = 0.17489547985835649
$ws.Range("J11").Value = 0.17489547985835646
$ws.Range("M11").Value = 0.01719833333333333
$ws.Range("N11").Value = 0.051595
$ws.Range("O11").Value = 0.0009451952024569835
$ws.Range("P11").Value = 0.0009451952024569837
$ws.Range("Q11").Value = 0.05473054853833332
$ws.Range("R11").Value = 0.49257493684500003
$ws.Range("S11").Value = 0.00016531036849353053
$ws.Range("T11").Value = 0.00016531036849353056

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 3.18231699999999984
$ws.Range("H12").Value = 9.54695099999999996
$ws.Range("I12").Value = 0.17489547985835649
$ws.Range("J12").Value = 0.17489547985835646
$ws.Range("M12").Value = 3.18231699999999984
$ws.Range("N12").Value = 9.54695099999999996
$ws.Range("O12").Value = 0.17489547985835649
$ws.Range("P12").Value = 0.17489547985835646
$ws.Range("Q12").Value = 10.12714148848899853
$ws.Range("R12").Value = 91.14427339640100456
$ws.Range("S12").Value = 0.03058842887488478
$ws.Range("T12").Value = 0.03058842887488477

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 3.18231699999999984
$ws.Range("H13").Value = 9.54695099999999996
$ws.Range("I13").Value = 0.17489547985835649
$ws.Range("J13").Value = 0.17489547985835646
$ws.Range("M13").Value = 14.17155133333332984
$ws.Range("N13").Value = 42.51465400000000017
$ws.Range("O13").Value = 0.77884769832190337
$ws.Range("P13").Value = 0.77884769832190348
$ws.Range("Q13").Value = 45.0983687244393181
$ws.Range("R13").Value = 405.88531851995401212
$ws.Range("S13").Value = 0.13621694193458575
$ws.Range("T13").Value = 0.13621694193458575

# Row 14: Resolving-Mac -> ECs
$ws.Range("G14").Value = 14.17155133333332984
$ws.Range("H14").Value = 42.51465400000000017
$ws.Range("I14").Value = 0.77884769832190337
$ws.Range("J14").Value = 0.77884769832190348
$ws.Range("M14").Value = 0.82446933333333339
$ws.Range("N14").Value = 2.47340800000000005
$ws.Range("O14").Value = 0.045311626617283135
$ws.Range("P14").Value = 0.04531162661728313
$ws.Range("Q14").Value = 11.68400948009244189
$ws.Range("R14").Value = 105.15608532083200544
$ws.Range("S14").Value = 0.03529085609809246
$ws.Range("T14").Value = 0.03529085609809246

# Row 15: Resolving-Mac -> FAPs
$ws.Range("G15").Value = 14.17155133333332984
$ws.Range("H15").Value = 42.51465400000000017
$ws.Range("I15").Value = 0.77884769832190337
$ws.Range("J15").Value = 0.77884769832190348
$ws.Range("M15").Value = 0.01719833333333333
$ws.Range("N15").Value = 0.051595
$ws.Range("O15").Value = 0.0009451952024569835
$ws.Range("P15").Value = 0.0009451952024569837
$ws.Range("Q15").Value = 0.24372706368111099
$ws.Range("R15").Value = 2.19354357312999992
$ws.Range("S15").Value = 0.0007361631078985271
$ws.Range("T15").Value = 0.0007361631078985273

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("G16").Value = 14.17155133333332984
$ws.Range("H16").Value = 42.51465400000000017
$ws.Range("I16").Value = 0.77884769832190337
$ws.Range("J16").Value = 0.77884769832190348
$ws.Range("M16").Value = 3.18231699999999984
$ws.Range("N16").Value = 9.54695099999999996
$ws.Range("O16").Value = 0.17489547985835649
$ws.Range("P16").Value = 0.17489547985835646
$ws.Range("Q16").Value = 45.0983687244393181
$ws.Range("R16").Value = 405.88531851995401212
$ws.Range("S16").Value = 0.13621694193458575
$ws.Range("T16").Value = 0.13621694193458575

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("G17").Value = 14.17155133333332984
$ws.Range("H17").Value = 42.51465400000000017
$ws.Range("I17").Value = 0.77884769832190337
$ws.Range("J17").Value = 0.77884769832190348
$ws.Range("M17").Value = 14.17155133333332984
$ws.Range("N17").Value = 42.51465400000000017
$ws.Range("O17").Value = 0.77884769832190337
$ws.Range("P17").Value = 0.77884769832190348
$ws.Range("Q17").Value = 200.83286719330166648
$ws.Range("R17").Value = 1807.49580473971604988
$ws.Range("S17").Value = 0.60660373718132665
$ws.Range("T17").Value = 0.60660373718132676
